# Generate Report for Handback
# - Update the "Ready for handoff" status text (shared across sheets) to
#   "Handback transform failed" for the 1fc7884c... file row.
# - Populate the "Error Detail" column (P) for that row on the zh-cn and
#   de-de sheets with the handback/handoff filename-mismatch message.
# - Widen column P on those sheets to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The status text "Ready for handoff" appears on Overview!E3/F3 and on
# zh-cn!C3 / de-de!C3 for the 1fc7884c... row - update every occurrence.
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Error Detail (column P) for the 1fc7884c... row on each locale sheet.
$wsZhCn.Range("P3").Value = "Handback file name: vtzwp42q.wkt is different with handoff file name: 1fc7884c-d2f1-42d4-bf0c-c08cf2aacb19.7e49bd0cf7cb197ebf836a02c251de3a65901b92.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: vtzwp42q.wkt is different with handoff file name: 1fc7884c-d2f1-42d4-bf0c-c08cf2aacb19.7e49bd0cf7cb197ebf836a02c251de3a65901b92.de-de."

# Widen column P to fit the new message text (stored column width of 40).
# ColumnWidth is in characters; Excel pads it internally, so 39.17 here
# round-trips to an XML-stored width of 40 (matches how column A's
# width="40" round-trips to a ColumnWidth of 39.17 in this workbook).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
